# Regenerate save_data column G ("K") values.
# The underlying data source switched from using "Strike#" to "K", and
# std/mean were regenerated, resulting in recalculated s_vals written
# into column G for each data row (rows 2-30).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2  = 4
    3  = 6
    4  = 7
    5  = 4
    6  = 5
    7  = 9
    8  = 6
    9  = 7
    10 = 5
    11 = 3
    12 = 8
    13 = 7
    14 = 6
    15 = 8
    16 = 13
    17 = 4
    18 = 8
    19 = 2
    20 = 2
    21 = 8
    22 = 7
    23 = 6
    24 = 7
    25 = 2
    26 = 9
    27 = 7
    28 = 4
    29 = 1
    30 = 1
}

foreach ($row in $values.Keys) {
    $ws.Range("G$row").Value = $values[$row]
}
